# "Really the final one"
# - Reassign several tasks from Ramiro Alvarex-Cruz/Alvarez-Cruz to
#   Hector Villalpando / Samuel Glover / the new team member Nicholas Greco.
# - Rename the last Documentation task and reassign it to Nicholas Greco.
# - Merge the Status (K) column data validation into a single contiguous range.
# - Update the saved view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reassign "Assigned to" (column C) for various tasks ---------------
$ws.Range("C5").Value  = "Hector Villalpando"
$ws.Range("C6").Value  = "Hector Villalpando"
$ws.Range("C7").Value  = "Samuel Glover"
$ws.Range("C9").Value  = "Hector Villalpando"
$ws.Range("C11").Value = "Samuel Glover"
$ws.Range("C14").Value = "Nicholas Greco"
$ws.Range("C18").Value = "Samuel Glover"
$ws.Range("C19").Value = "Hector Villalpando"
$ws.Range("C20").Value = "Nicholas Greco"
$ws.Range("C40").Value = "Nicholas Greco"
$ws.Range("C43").Value = "Hector Villalpando"
$ws.Range("C46").Value = "Nicholas Greco"

# --- Rename the final "User Story #9: Documentation" task --------------
$ws.Range("A46").Value = "Document sprint reviews and backlogs"

# --- Merge the two Status dropdown validation ranges into one ----------
$ws.Range("K4:K46").Validation.Delete()
$ws.Range("K4:K46").Validation.Add(3, 1, 1, '"To Do,Doing,Done"')
$ws.Range("K4:K46").Validation.IgnoreBlank = $true
$ws.Range("K4:K46").Validation.InCellDropdown = $true
$ws.Range("K4:K46").Validation.ShowInput = $false
$ws.Range("K4:K46").Validation.ShowError = $true

# --- Update the view / selection state ----------------------------------
$ws.Range("C37").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
